$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new column at D, shifting City..Address (old D:H) to E:I
$ws.Columns.Item(4).Insert()

# Give the new column D a sensible width similar to its neighbours
$ws.Columns.Item(4).ColumnWidth = 10.3

# Update existing donor record values
$ws.Range("B2").Value = 1234567890
$ws.Range("C2").Value = "Hyderabad"

# New "Area" column header and value
$ws.Range("D1").Value = "Area"
$ws.Range("D2").Value = "kukatlapally"

# Update the active selection on the sheet
$ws.Range("G13").Select()
